$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45951
$ws.Range("B2").Value = 35.95
$ws.Range("C2").Value = 24.16
$ws.Range("D2").Value = 13.08
$ws.Range("E2").Value = 11.05
$ws.Range("F2").Value = 13.98
$ws.Range("G2").Value = 18.04
$ws.Range("H2").Value = 35.27
$ws.Range("I2").Value = 62.97
$ws.Range("J2").Value = 93.41
$ws.Range("K2").Value = 78.52
$ws.Range("L2").Value = 41.49
$ws.Range("M2").Value = 19.19
$ws.Range("N2").Value = 5.86
$ws.Range("O2").Value = 3.52
$ws.Range("P2").Value = 3.26
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 3.39
$ws.Range("S2").Value = 13.88
$ws.Range("T2").Value = 48.03
$ws.Range("U2").Value = 78.31
$ws.Range("V2").Value = 101.19
$ws.Range("W2").Value = 82.97
$ws.Range("X2").Value = 78.82
$ws.Range("Y2").Value = 52.95
$ws.Range("Z2").Value = 38.38
$ws.Range("AB2").Value = 78.98
$ws.Range("AD2").Value = 92.08
$ws.Range("AF2").Value = 85.96
$ws.Range("AG2").Value = "0h-17h"
